# Add a new "chPanel" row to the observations summary sheet, inserted
# right above the existing "chTest" row, and rename the "chTest" row's
# Name cell to include ".chemistryResults".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (pushes chTest and everything below it
# down by one row).
$ws.Rows.Item(2).Insert()

# Copy formatting (style/borders) from the row below (the old row 2,
# now row 3) into the newly-inserted row 2 so it keeps the same table
# styling as the rest of the data rows.
$ws.Range("A3:K3").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

# Populate the new chPanel row.
$ws.Range("A2").Value = "VA.MHV.PHR.chPanel"
$ws.Range("B2").Value = "VA MHV PHR CH labTests.orderedTestCode"
$ws.Range("C2").Value = "Observation Category Codes#laboratory"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G2").Value = "dateTime, Period, Timing, instant"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "prohibited"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

# Update the chTest row's Name cell (now row 3 after the insert).
$ws.Range("B3").Value = "VA MHV PHR CH labTest.chemistryResults"
